$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: 18-Dec-2019, "design the home page" (new text), D = "Design the Home page"
# Copy A10's format first (keeps the existing date style index) then overwrite the value.
$ws.Range("A10").Copy($ws.Range("A11")) | Out-Null
$ws.Range("A11").Value = 43817
$ws.Range("B11").Value = "design the home page"
$ws.Range("D11").Value = "Design the Home page"

# Row 12: 19-Dec-2019, "Edit the document", "DFD diagram"
$ws.Range("A10").Copy($ws.Range("A12")) | Out-Null
$ws.Range("A12").Value = 43818
$ws.Range("B12").Value = "Edit the document"
$ws.Range("C12").Value = "DFD diagram"

# Update the view: drop the frozen top-left cell and move the active selection
$ws.Range("C13").Select() | Out-Null

Write-Output "done"
